# Casting : ajout bouton refresh
# Applies the textual / formatting corrections captured in the diff:
#  - straight double quotes -> French guillemets (« »), on three comments
#  - "Le 7ème" -> "Le 7" + superscript "e" + " cercle..." (rich text run)
#  - row 54 height bump 13.8 -> 14.2
#  - selection moved to A71 (scroll/selection bookkeeping)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A25: "de manière naturelle" -> « de manière naturelle » -----------
$ws.Range("A25").Value = [char]0x00AB + " de manière naturelle " + [char]0x00BB

# --- A52: "Y a une petite mise en scène" ... -> « Y a une petite mise en scène. » ...
$ws.Range("A52").Value = [char]0x00AB + " Y a une petite mise en scène. " + [char]0x00BB + " Ah ouais vraiment toute petite !"

# --- A124: petit "putain" -> petit « putain » ---------------------------
# (keeps the existing non-breaking space before the colon, per the
# document's French-typography convention used throughout the sheet)
$ws.Range("A124").Value = "Atouts" + [char]0x00A0 + ": il fait très bien la goutte, petit " + [char]0x00AB + " putain " + [char]0x00BB + " quand il oublie son texte"

# --- A54: "Le 7ème cercle..." -> "Le 7" + superscript "e" + " cercle..." -
$cell54 = $ws.Range("A54")
$cell54.Value = "Le 7e cercle de l'enfer est constitué de diffusion en boucle de scènes de théâtre joués par ces protagonistes"

$runSup = $cell54.Characters(5, 1)
$runSup.Font.Name = "Cambria"
$runSup.Font.Size = 11
$runSup.Font.Superscript = $true
$runSup.Font.ColorIndex = -4105

$runBefore = $cell54.Characters(1, 4)
$runBefore.Font.Name = "Cambria"
$runBefore.Font.Size = 11
$runBefore.Font.ColorIndex = -4105

$runAfter = $cell54.Characters(6, 104)
$runAfter.Font.Name = "Cambria"
$runAfter.Font.Size = 11
$runAfter.Font.ColorIndex = -4105

# --- row 54 height 13.8 -> 14.2 -----------------------------------------
$ws.Rows.Item(54).RowHeight = 14.2

# --- selection / scroll bookkeeping (view-state) -------------------------
$ws.Activate() | Out-Null
$ws.Range("A71").Select() | Out-Null
